$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("G1-L1")
$ws.Range("C3").Value = ""
$ws.Range("C3").Interior.Pattern = -4142
$ws.Range("G3").Value = ""
$ws.Range("G3").Interior.Pattern = -4142
$ws.Range("G5").Value = "[DPR110] communication`n(CM)`nProf: Dieynaba`nSalle: 101"
$ws.Range("G5").Interior.Color = 10086143
$ws.Range("G7").Value = "[DPR110] communication`n(CM)`nProf: Dieynaba`nSalle: 101"
$ws.Range("G7").Interior.Color = 10086143

$ws = $wb.Worksheets.Item("G2-L1")
$ws.Range("C5").Value = ""
$ws.Range("C5").Interior.Pattern = -4142
$ws.Range("E6").Value = ""
$ws.Range("E6").Interior.Pattern = -4142
$ws.Range("B7").Value = "[DPR110] communication`n(CM)`nProf: Dieynaba`nSalle: 101"
$ws.Range("B7").Interior.Color = 10086143
$ws.Range("D7").Value = "[DPR110] communication`n(CM)`nProf: Dieynaba`nSalle: 101"
$ws.Range("D7").Interior.Color = 10086143

$ws = $wb.Worksheets.Item("G3-L1")
$ws.Range("G3").Value = "[DPR110] communication`n(CM)`nProf: Dieynaba`nSalle: 101"
$ws.Range("G3").Interior.Color = 10086143
$ws.Range("C7").Value = "[DPR110] communication`n(CM)`nProf: Dieynaba`nSalle: 101"
$ws.Range("C7").Interior.Color = 10086143
$ws.Range("D7").Value = ""
$ws.Range("D7").Interior.Pattern = -4142
$ws.Range("F7").Value = ""
$ws.Range("F7").Interior.Pattern = -4142

$ws = $wb.Worksheets.Item("G4-L1")
$ws.Range("D5").Value = ""
$ws.Range("D5").Interior.Pattern = -4142
$ws.Range("E6").Value = "[DPR110] communication`n(CM)`nProf: Dieynaba`nSalle: 101"
$ws.Range("E6").Interior.Color = 10086143

$ws = $wb.Worksheets.Item("RSS-L2")
$ws.Range("B3").Value = "[RSS310] Reseaux Mobile`n(CM)`nProf: Aloun`nSalle: 101"
$ws.Range("B3").Interior.Color = 10086143
$ws.Range("C3").Value = "[DAS311] RO`n(TD) - TD1`nProf: abderrahmane`nSalle: 101 /// [RSS310] Reseaux Mobile`n(TD) - TD2`nProf: Aloun`nSalle: 102"
$ws.Range("C3").Interior.Color = 15189940
$ws.Range("D3").Value = "[DPR310] Communication`n(CM)`nProf: Dieynaba`nSalle: 101"
$ws.Range("D3").Interior.Color = 10086143
$ws.Range("E3").Value = "[PAV312] Projet Integrateur`n(CM)`nProf: Encadreur`nSalle: 101"
$ws.Range("E3").Interior.Color = 10086143
$ws.Range("F3").Value = "[RSS320] Introduction a la securite`n(TP) - TD1`nProf: Tourad`nSalle: 102 /// [DAS311] RO`n(TP) - TD2`nProf: abderrahmane`nSalle: 103"
$ws.Range("F3").Interior.Color = 11854022
$ws.Range("G3").Value = "[DAS311] RO`n(TD) - TD1`nProf: abderrahmane`nSalle: 102 /// [RSS320] Introduction a la securite`n(TD) - TD2`nProf: Tourad`nSalle: 103"
$ws.Range("G3").Interior.Color = 15189940
$ws.Range("B4").Value = "[DPR310] Communication`n(CM Online)`nProf: Dieynaba`nSalle: En ligne"
$ws.Range("B4").Interior.Pattern = -4142
$ws.Range("C4").Value = "[RSS311] Administration reseaux`n(CM)`nProf: Aloun`nSalle: 101"
$ws.Range("C4").Interior.Color = 10086143
$ws.Range("D4").Value = "[PAV311] SD & Comp.Algo`n(CM)`nProf: Meyara`nSalle: 201"
$ws.Range("D4").Interior.Color = 10086143
$ws.Range("F4").Value = "[RSS321] BD & CSI`n(CM)`nProf: Med Lemine`nSalle: 101"
$ws.Range("F4").Interior.Color = 10086143
$ws.Range("G4").Value = "[PAV312] Projet Integrateur`n(CM)`nProf: Encadreur`nSalle: 101"
$ws.Range("G4").Interior.Color = 10086143
$ws.Range("C5").Value = "[DAS310] Maching Learning`n(CM Online)`nProf: Louly`nSalle: En ligne"
$ws.Range("C5").Interior.Pattern = -4142
$ws.Range("D5").Value = "[RSS310] Reseaux Mobile`n(TD) - TD1`nProf: Aloun`nSalle: 101 /// [DAS311] RO`n(TD) - TD2`nProf: abderrahmane`nSalle: 102"
$ws.Range("D5").Interior.Color = 15189940
$ws.Range("E5").Value = "[RSS310] Reseaux Mobile`n(TD) - TD1`nProf: Aloun`nSalle: 101 /// [DAS311] RO`n(TD) - TD2`nProf: abderrahmane`nSalle: 102"
$ws.Range("E5").Interior.Color = 15189940
$ws.Range("G5").Value = "[RSS311] Administration reseaux`n(TP) - TD1`nProf: Aloun`nSalle: 102 /// [RSS320] Introduction a la securite`n(TP) - TD2`nProf: Tourad`nSalle: 103"
$ws.Range("G5").Interior.Color = 11854022
$ws.Range("B6").Value = "[RSS321] BD & CSI`n(TD Online) - TD1`nProf: Med Lemine`nSalle: En ligne /// [DAS311] RO`n(TD) - TD2`nProf: abderrahmane`nSalle: 101"
$ws.Range("B6").Interior.Color = 15189940
$ws.Range("C6").Value = "[RSS320] Introduction a la securite`n(TD) - TD1`nProf: Tourad`nSalle: 101 /// [RSS310] Reseaux Mobile`n(TD) - TD2`nProf: Aloun`nSalle: 102"
$ws.Range("C6").Interior.Color = 15189940
$ws.Range("D6").Value = "[DAS311] RO`n(CM)`nProf: Cheikh`nSalle: 101"
$ws.Range("D6").Interior.Color = 10086143
$ws.Range("E6").Value = "[RSS321] BD & CSI`n(CM)`nProf: Med Lemine`nSalle: 201"
$ws.Range("E6").Interior.Color = 10086143
$ws.Range("G6").Value = "[PAV310] POO JAVA`n(CM)`nProf: Esseyssah`nSalle: 101"
$ws.Range("G6").Interior.Color = 10086143
$ws.Range("B7").Value = "[DPR313] Gestion d'enterprise`n(CM)`nProf: El Bennany`nSalle: 201"
$ws.Range("B7").Interior.Color = 10086143
$ws.Range("C7").Value = "[RSS320] Introduction a la securite`n(TP) - TD1`nProf: Tourad`nSalle: 102 /// [RSS311] Administration reseaux`n(TP) - TD2`nProf: Aloun`nSalle: 103"
$ws.Range("C7").Interior.Color = 11854022
$ws.Range("D7").Value = "[DAS311] RO`n(TD) - TD1`nProf: abderrahmane`nSalle: 102 /// [RSS321] BD & CSI`n(TD Online) - TD2`nProf: Med Lemine`nSalle: En ligne"
$ws.Range("D7").Interior.Color = 15189940
$ws.Range("E7").Value = "[DAS311] RO`n(TP) - TD1`nProf: abderrahmane`nSalle: 102 /// [RSS320] Introduction a la securite`n(TP) - TD2`nProf: Tourad`nSalle: 103"
$ws.Range("E7").Interior.Color = 11854022
$ws.Range("F7").Value = "[DPR310] Communication`n(CM)`nProf: Dieynaba`nSalle: 101"
$ws.Range("F7").Interior.Color = 10086143
$ws.Range("G7").Value = "[DAS310] Maching Learning`n(CM)`nProf: Louly`nSalle: 201"
$ws.Range("G7").Interior.Color = 10086143
